# Refactor cargar_datos: mapeo directo de columnas Excel, simplificado y robusto
#
# Adds the "Temperatura ambiente (°C)" (column T) and "Humedad Relativa (%)"
# (column U) readings for each data row (rows 2-28) on Hoja1, mirroring the
# direct column-to-column mapping used by the refactored data loader.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$temperaturas = @{
    2  = 20
    3  = 10
    4  = 30
    5  = 25
    6  = 16
    7  = 17
    8  = 18
    9  = 19
    10 = 20
    11 = 21
    12 = 22
    13 = 23
    14 = 24
    15 = 25
    16 = 26
    17 = 27
    18 = 28
    19 = 29
    20 = 30
    21 = 31
    22 = 32
    23 = 33
    24 = 34
    25 = 35
    26 = 36
    27 = 37
    28 = 38
}

$humedades = @{
    2  = 20
    3  = 21
    4  = 22
    5  = 23
    6  = 24
    7  = 25
    8  = 26
    9  = 27
    10 = 28
    11 = 29
    12 = 30
    13 = 31
    14 = 32
    15 = 33
    16 = 34
    17 = 35
    18 = 36
    19 = 37
    20 = 38
    21 = 39
    22 = 40
    23 = 41
    24 = 42
    25 = 43
    26 = 44
    27 = 45
    28 = 46
}

foreach ($fila in 2..28) {
    $ws.Range("T$fila").Value = $temperaturas[$fila]
    $ws.Range("U$fila").Value = $humedades[$fila]
}

# Restore the active sheet view / selection as left by the editing session.
$ws.Range("T21").Select()
